$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This workbook tracks several "everyday" data tables, one table per sheet,
# each sorted newest-date-first starting at row 2 (row 1 is the header).
# The edit adds one new day (2021/12/14) at the top of every table: existing
# rows 2..N shift down to 3..N+1, and the brand-new data is written into the
# now-empty row 2.
# ---------------------------------------------------------------------------

function Shift-DownAndInsert($SheetIndex, $LastCol, $LastDataRow, $NewRow) {
    $ws = $wb.Worksheets.Item($SheetIndex)

    # Read the existing data block (row 2 .. LastDataRow) in one shot, then
    # write it back one row lower so rows 2..LastDataRow become 3..LastDataRow+1.
    $srcRange = $ws.Range("A2:" + $LastCol + $LastDataRow)
    $vals = $srcRange.Value()

    $dstRange = $ws.Range("A3:" + $LastCol + ($LastDataRow + 1))
    $dstRange.Value = $vals

    # Write the brand-new first row.
    $lastColNum = $NewRow.Length
    for ($c = 1; $c -le $lastColNum; $c++) {
        $ws.Cells.Item(2, $c).Value = $NewRow[$c - 1]
    }
}

# Sheet 1: 台指期換倉成本計算 (A1:F20 -> A1:F21), data rows 2..19 -> 3..20
Shift-DownAndInsert 1 "F" 19 @("日期：2021/12/14", "202201", 17543, 59283, 367332877, 17635)

# Sheet 2: 散戶多空力道 (A1:B35 -> A1:B36), data rows 2..35 -> 3..36
Shift-DownAndInsert 2 "B" 35 @("日期：2021/12/14", 0.06)

# Sheet 3: 三大法人買賣金額 (A1:C35 -> A1:C36), data rows 2..35 -> 3..36
Shift-DownAndInsert 3 "C" 35 @("110年12月14日", -156, -12.23)

# Sheet 4: 大盤多空點位 (A1:B34 -> A1:B35), data rows 2..34 -> 3..35
Shift-DownAndInsert 4 "B" 34 @("110年12月14日", 17635.97)

# Sheet 5: 期貨大額交易人未沖銷部位 (A1:N33 -> A1:N34), data rows 2..33 -> 3..34
Shift-DownAndInsert 5 "N" 33 @("2021/12/14", 51761, 59284, 800, 1430, 26018, 52546, -995, -217, -26528, -778, 1795, 1647, 148)
